# Insert a new data row at row 229 (pushing existing rows 229-272 down to
# 230-273) and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 229, shifting rows 229:272 down to 230:273
$ws.Rows.Item(229).Insert()

$newRow = 229
$ws.Cells.Item($newRow, 1).Value = 6
$ws.Cells.Item($newRow, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item($newRow, 3).Value = "Metropolitana"
$ws.Cells.Item($newRow, 4).Value = 44476
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat
$ws.Cells.Item($newRow, 5).Value = 13
$ws.Cells.Item($newRow, 6).Value = 100112030
$ws.Cells.Item($newRow, 7).Value = "Poroto granado"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 150
$ws.Cells.Item($newRow, 11).Value = 75000
$ws.Cells.Item($newRow, 12).Value = 80000
$ws.Cells.Item($newRow, 13).Value = 77667
$ws.Cells.Item($newRow, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item($newRow, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($newRow, 16).Value = 3107
$ws.Cells.Item($newRow, 17).Value = 25
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
